# This workbook contains a weekly price-report dataset on Sheet1.
# A new week's record is inserted at row 413 (just after the header block
# of earlier, unrelated rows), pushing the existing rows 413-469 down to
# 414-470 and extending the used range to A1:R470.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 413; this shifts rows 413:469 -> 414:470
$ws.Rows.Item(413).Insert()

# Populate the newly inserted row 413 with the new weekly record
$ws.Cells.Item(413, 1).Value = 4
$ws.Cells.Item(413, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(413, 3).Value = "Los Lagos"
$ws.Cells.Item(413, 4).Value = 45124
$ws.Cells.Item(413, 5).Value = 10
$ws.Cells.Item(413, 6).Value = 100112043
$ws.Cells.Item(413, 7).Value = "Pepino ensalada"
$ws.Cells.Item(413, 8).Value = "Sin especificar"
$ws.Cells.Item(413, 9).Value = "Primera"
$ws.Cells.Item(413, 10).Value = 200
$ws.Cells.Item(413, 11).Value = 17000
$ws.Cells.Item(413, 12).Value = 17000
$ws.Cells.Item(413, 13).Value = 17000
$ws.Cells.Item(413, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(413, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(413, 16).Value = 283
$ws.Cells.Item(413, 17).Value = 60
$ws.Cells.Item(413, 18).Value = "Hortaliza"
